$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24 has no cell in column A yet - copy formatting from an existing
# cell in the row (B24) so the new A24 cell picks up the same style (s="2").
$ws.Range("B24").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A24").Value = "Invitation mails for shusha chess"
$ws.Range("B24").Value = "Shusha"
$ws.Range("C24").Value = "High"
$ws.Range("D24").Value = "Medium"
$ws.Range("E24").Value = "Low"
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = "Not Started"
$ws.Range("H24").Value = 44789

$ws.Range("H25").Select()
